$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5231.25
$ws.Range("I40").Value = 11400
$ws.Range("J40").Value = 3175
$ws.Range("K40").Value = 11400
$ws.Range("L40").Value = 3175
$ws.Range("M40").Value = -11225
$ws.Range("N40").Value = -3525
$ws.Range("H100").Value = 2270.9546
$ws.Range("I100").Value = 1519
$ws.Range("J100").Value = 3022.9092
$ws.Range("K100").Value = 1519
$ws.Range("L100").Value = 3022.9092
$ws.Range("M100").Value = -978
$ws.Range("N100").Value = -4104.9092
$ws.Range("H107").Value = 335.2
$ws.Range("I107").Value = 218.17647
$ws.Range("K107").Value = 218.17647
$ws.Range("M107").Value = 1701.82353
$ws.Range("H129").Value = 8330.333000000001
$ws.Range("J129").Value = 8899.036
$ws.Range("L129").Value = 26697.108
$ws.Range("N129").Value = -36697.108
$ws.Range("H137").Value = 1070.2222
$ws.Range("I137").Value = 835.11536
$ws.Range("J137").Value = 1681.5
$ws.Range("K137").Value = 2505.34608
$ws.Range("L137").Value = 5044.5
$ws.Range("M137").Value = 44.65391999999974
$ws.Range("N137").Value = -10144.5
$ws.Range("H138").Value = 4265.015
$ws.Range("I138").Value = 2023.2693
$ws.Range("J138").Value = 5722.15
$ws.Range("K138").Value = 6069.8079
$ws.Range("L138").Value = 17166.45
$ws.Range("M138").Value = -929.8078999999998
$ws.Range("N138").Value = -27446.45

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3348.6667
$ws.Range("I74").Value = 3235.9048
$ws.Range("K74").Value = 3235.9048
$ws.Range("M74").Value = -2361.9048
$ws.Range("H77").Value = 3348.6667
$ws.Range("I77").Value = 3235.9048
$ws.Range("K77").Value = 16179.524
$ws.Range("M77").Value = -11811.524
$ws.Range("H97").Value = 2399.12
$ws.Range("I97").Value = 2691.4211
$ws.Range("J97").Value = 1473.5
$ws.Range("K97").Value = 2691.4211
$ws.Range("L97").Value = 1473.5
$ws.Range("M97").Value = -2195.4211
$ws.Range("N97").Value = -2465.5
$ws.Range("H122").Value = 2787.804
$ws.Range("I122").Value = 2464.6047
$ws.Range("J122").Value = 4525
$ws.Range("K122").Value = 7393.8141
$ws.Range("L122").Value = 13575
$ws.Range("M122").Value = -4943.8141
$ws.Range("N122").Value = -18475

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3560.65
$ws.Range("I20").Value = 1564.7858
$ws.Range("J20").Value = 8217.666999999999
$ws.Range("K20").Value = 1564.7858
$ws.Range("L20").Value = 8217.666999999999
$ws.Range("M20").Value = -1317.7858
$ws.Range("N20").Value = -8711.666999999999
$ws.Range("H86").Value = 4001931.2
$ws.Range("I86").Value = 5265102
$ws.Range("J86").Value = 1890.5834
$ws.Range("K86").Value = 5265102
$ws.Range("L86").Value = 1890.5834
$ws.Range("M86").Value = -5263979
$ws.Range("N86").Value = -4136.5834
$ws.Range("H89").Value = 4001931.2
$ws.Range("I89").Value = 5265102
$ws.Range("J89").Value = 1890.5834
$ws.Range("K89").Value = 26325510
$ws.Range("L89").Value = 9452.916999999999
$ws.Range("M89").Value = -26319894
$ws.Range("N89").Value = -20684.917

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2404.3157
$ws.Range("I31").Value = 1282.6086
$ws.Range("J31").Value = 4124.2666
$ws.Range("K31").Value = 1282.6086
$ws.Range("L31").Value = 4124.2666
$ws.Range("M31").Value = -987.6086
$ws.Range("N31").Value = -4714.2666
$ws.Range("H34").Value = 2404.3157
$ws.Range("I34").Value = 1282.6086
$ws.Range("J34").Value = 4124.2666
$ws.Range("K34").Value = 1282.6086
$ws.Range("L34").Value = 4124.2666
$ws.Range("M34").Value = -1080.6086
$ws.Range("N34").Value = -4528.2666
$ws.Range("H60").Value = 17660
$ws.Range("J60").Value = 17660
$ws.Range("L60").Value = 17660
$ws.Range("N60").Value = -18682
$ws.Range("H94").Value = 1040.2142
$ws.Range("I94").Value = 977
$ws.Range("J94").Value = 1103.4286
$ws.Range("K94").Value = 977
$ws.Range("L94").Value = 1103.4286
$ws.Range("M94").Value = -526
$ws.Range("N94").Value = -2005.4286
$ws.Range("H99").Value = 10258.4
$ws.Range("I99").Value = 19156
$ws.Range("J99").Value = 4326.6665
$ws.Range("K99").Value = 19156
$ws.Range("L99").Value = 4326.6665
$ws.Range("M99").Value = -17658
$ws.Range("N99").Value = -7322.6665
$ws.Range("H126").Value = 10258.4
$ws.Range("I126").Value = 19156
$ws.Range("J126").Value = 4326.6665
$ws.Range("K126").Value = 57468
$ws.Range("L126").Value = 12979.9995
$ws.Range("M126").Value = -54998
$ws.Range("N126").Value = -17919.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1423.8
$ws.Range("I7").Value = 406.66666
$ws.Range("J7").Value = 2949.5
$ws.Range("K7").Value = 1219.99998
$ws.Range("L7").Value = 8848.5
$ws.Range("M7").Value = -1107.99998
$ws.Range("N7").Value = -9072.5
$ws.Range("H17").Value = 476.8
$ws.Range("J17").Value = 496
$ws.Range("L17").Value = 1488
$ws.Range("N17").Value = -1826
$ws.Range("H34").Value = 562.6316
$ws.Range("I34").Value = 67.27273
$ws.Range("J34").Value = 1243.75
$ws.Range("K34").Value = 201.81819
$ws.Range("L34").Value = 3731.25
$ws.Range("M34").Value = -117.81819
$ws.Range("N34").Value = -3899.25
$ws.Range("H39").Value = 2835.9285
$ws.Range("J39").Value = 3100
$ws.Range("L39").Value = 9300
$ws.Range("N39").Value = -9888
$ws.Range("H55").Value = 2093.0476
$ws.Range("J55").Value = 2162.5
$ws.Range("L55").Value = 6487.5
$ws.Range("N55").Value = -6841.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2477.7778
$ws.Range("I80").Value = 2450
$ws.Range("K80").Value = 2450
$ws.Range("M80").Value = -1452
$ws.Range("H83").Value = 2477.7778
$ws.Range("I83").Value = 2450
$ws.Range("K83").Value = 12250
$ws.Range("M83").Value = -7258
$ws.Range("H122").Value = 2269.1667
$ws.Range("I122").Value = 3163.2
$ws.Range("J122").Value = 1630.5714
$ws.Range("K122").Value = 9489.599999999999
$ws.Range("L122").Value = 4891.7142
$ws.Range("M122").Value = -7039.599999999999
$ws.Range("N122").Value = -9791.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 12033
$ws.Range("J38").Value = 12033
$ws.Range("L38").Value = 12033
$ws.Range("N38").Value = -12853
$ws.Range("H40").Value = 2142
$ws.Range("I40").Value = 1868.6666
$ws.Range("K40").Value = 1868.6666
$ws.Range("M40").Value = -1732.6666
$ws.Range("H93").Value = 2032.4117
$ws.Range("I93").Value = 2264.7144
$ws.Range("J93").Value = 1869.8
$ws.Range("K93").Value = 2264.7144
$ws.Range("L93").Value = 1869.8
$ws.Range("M93").Value = -1016.7144
$ws.Range("N93").Value = -4365.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5139.3335
$ws.Range("I62").Value = 4950
$ws.Range("J62").Value = 5355.7144
$ws.Range("K62").Value = 4950
$ws.Range("L62").Value = 5355.7144
$ws.Range("M62").Value = -4326
$ws.Range("N62").Value = -6603.7144
$ws.Range("H65").Value = 5139.3335
$ws.Range("I65").Value = 4950
$ws.Range("J65").Value = 5355.7144
$ws.Range("K65").Value = 24750
$ws.Range("L65").Value = 26778.572
$ws.Range("M65").Value = -21630
$ws.Range("N65").Value = -33018.572
$ws.Range("H126").Value = 1277.5454
$ws.Range("I126").Value = 1899
$ws.Range("J126").Value = 1139.4445
$ws.Range("K126").Value = 5697
$ws.Range("L126").Value = 3418.3335
$ws.Range("M126").Value = -3227
$ws.Range("N126").Value = -8358.333500000001
